# Add the "budget" model (Model_ID=4) rows to the schema table.
# Mirrors: Added budget model (Model_ID=4) with 8 fields to the schema sheet,
# expanding the "schema" table from A1:F46 to A1:F54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(401, 4, "budget", "Account_id",     "Account Number", $true),
  @(402, 4, "budget", "Account_Name",   "Account Name",   $true),
  @(403, 4, "budget", "Month",          "Month",          $true),
  @(404, 4, "budget", "Entity",         "Entity",         $true),
  @(405, 4, "budget", "Classification", "Classification", $true),
  @(406, 4, "budget", "EBITA",          "EBITA",          $true),
  @(407, 4, "budget", "id",             "id",             $true),
  @(408, 4, "budget", "Amount",         "Amount",         $true)
)

$r = 47
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# Grow the "schema" table (and its autofilter) to cover the new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F54"))

# Match the saved selection/view state (new rows highlighted).
$ws.Range("F47:F54").Select()
